# [2022-09-30] - Monkeypox update (code)
# Weekly RIVM-style excess-mortality-by-province refresh:
#  - minor upward revisions to several already-reported weeks (rows 123-142)
#  - two brand-new weeks appended: 2022 week 37 (row 143) and week 38 (row 144)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised raw counts on previously-filed weeks ---
$ws.Range("U123").Value = 394
$ws.Range("W125").Value = 515
$ws.Range("X126").Value = 573
$ws.Range("X129").Value = 593
$ws.Range("S132").Value = 209
$ws.Range("X132").Value = 573
$ws.Range("Z132").Value = 500
$ws.Range("AA132").Value = 260
$ws.Range("X133").Value = 605
$ws.Range("U135").Value = 403
$ws.Range("W135").Value = 509
$ws.Range("X135").Value = 638
$ws.Range("T136").Value = 64
$ws.Range("W136").Value = 511
$ws.Range("X136").Value = 599
$ws.Range("Z136").Value = 476
$ws.Range("T137").Value = 40
$ws.Range("V137").Value = 198
$ws.Range("W137").Value = 453
$ws.Range("X137").Value = 612
$ws.Range("Z137").Value = 444
$ws.Range("AA137").Value = 218
$ws.Range("W138").Value = 446
$ws.Range("X138").Value = 592
$ws.Range("Y138").Value = 79
$ws.Range("U139").Value = 381
$ws.Range("V139").Value = 210
$ws.Range("W139").Value = 500
$ws.Range("X139").Value = 603
$ws.Range("Z139").Value = 474
$ws.Range("AA139").Value = 205
$ws.Range("R140").Value = 109
$ws.Range("U140").Value = 357
$ws.Range("V140").Value = 207
$ws.Range("W140").Value = 441
$ws.Range("X140").Value = 594
$ws.Range("Z140").Value = 488
$ws.Range("U141").Value = 393
$ws.Range("V141").Value = 188
$ws.Range("W141").Value = 407
$ws.Range("X141").Value = 595
$ws.Range("Y141").Value = 72
$ws.Range("Z141").Value = 476
$ws.Range("P142").Value = 122
$ws.Range("Q142").Value = 122
$ws.Range("R142").Value = 91
$ws.Range("S142").Value = 201
$ws.Range("T142").Value = 44
$ws.Range("V142").Value = 186
$ws.Range("W142").Value = 451
$ws.Range("X142").Value = 592
$ws.Range("Y142").Value = 79
$ws.Range("Z142").Value = 449
$ws.Range("AA142").Value = 209

# --- New week: 2022 week 37 (row 143) ---
$ws.Range("N143").Value = 2022
$ws.Range("O143").Value = 37
$ws.Range("P143").Value = 110
$ws.Range("Q143").Value = 109
$ws.Range("R143").Value = 85
$ws.Range("S143").Value = 191
$ws.Range("T143").Value = 46
$ws.Range("U143").Value = 368
$ws.Range("V143").Value = 169
$ws.Range("W143").Value = 399
$ws.Range("X143").Value = 596
$ws.Range("Y143").Value = 60
$ws.Range("Z143").Value = 417
$ws.Range("AA143").Value = 221
$ws.Range("AC143").Value = 2022
$ws.Range("AD143").Value = 37

$ws.Range("AE143").Formula = "=ROUND((P143-B143)/B143*100,2)"
$ws.Range("AF143").Formula = "=ROUND((Q143-C143)/C143*100,2)"
$ws.Range("AG143").Formula = "=ROUND((R143-D143)/D143*100,2)"
$ws.Range("AH143").Formula = "=ROUND((S143-E143)/E143*100,2)"
$ws.Range("AI143").Formula = "=ROUND((T143-F143)/F143*100,2)"
$ws.Range("AJ143").Formula = "=ROUND((U143-G143)/G143*100,2)"
$ws.Range("AK143").Formula = "=ROUND((V143-H143)/H143*100,2)"
$ws.Range("AL143").Formula = "=ROUND((W143-I143)/I143*100,2)"
$ws.Range("AM143").Formula = "=ROUND((X143-J143)/J143*100,2)"
$ws.Range("AN143").Formula = "=ROUND((Y143-K143)/K143*100,2)"
$ws.Range("AO143").Formula = "=ROUND((Z143-L143)/L143*100,2)"
$ws.Range("AP143").Formula = "=ROUND((AA143-M143)/M143*100,2)"

# --- New week: 2022 week 38 (row 144) ---
$ws.Range("N144").Value = 2022
$ws.Range("O144").Value = 38
$ws.Range("P144").Value = 101
$ws.Range("Q144").Value = 140
$ws.Range("R144").Value = 109
$ws.Range("S144").Value = 196
$ws.Range("T144").Value = 46
$ws.Range("U144").Value = 362
$ws.Range("V144").Value = 213
$ws.Range("W144").Value = 487
$ws.Range("X144").Value = 575
$ws.Range("Y144").Value = 69
$ws.Range("Z144").Value = 423
$ws.Range("AA144").Value = 221
$ws.Range("AC144").Value = 2022
$ws.Range("AD144").Value = 38

$ws.Range("AE144").Formula = "=ROUND((P144-B144)/B144*100,2)"
$ws.Range("AF144").Formula = "=ROUND((Q144-C144)/C144*100,2)"
$ws.Range("AG144").Formula = "=ROUND((R144-D144)/D144*100,2)"
$ws.Range("AH144").Formula = "=ROUND((S144-E144)/E144*100,2)"
$ws.Range("AI144").Formula = "=ROUND((T144-F144)/F144*100,2)"
$ws.Range("AJ144").Formula = "=ROUND((U144-G144)/G144*100,2)"
$ws.Range("AK144").Formula = "=ROUND((V144-H144)/H144*100,2)"
$ws.Range("AL144").Formula = "=ROUND((W144-I144)/I144*100,2)"
$ws.Range("AM144").Formula = "=ROUND((X144-J144)/J144*100,2)"
$ws.Range("AN144").Formula = "=ROUND((Y144-K144)/K144*100,2)"
$ws.Range("AO144").Formula = "=ROUND((Z144-L144)/L144*100,2)"
$ws.Range("AP144").Formula = "=ROUND((AA144-M144)/M144*100,2)"

# --- Scroll / selection state, matching the author's saved view ---
$ws.Application.Goto($ws.Range("A127"), $false)
$ws.Range("AL144").Select()
